$wb = $excel.ActiveWorkbook

# --- Data fixes on "ShowDateRange" sheet (formerly test data with placeholder y/n flags) ---
$wsShowDateRange = $wb.Worksheets.Item("ShowDateRange")
# H4 was "y" -> now "Gar" (entered with a leading apostrophe so Excel stores it as
# quote-prefixed text, matching the quotePrefix style seen in the target file)
$wsShowDateRange.Range("H4").Value = "'Gar"
# H6 was "n" -> now "ET"
$wsShowDateRange.Range("H6").Value = "ET"

# --- Data fixes on "AdvanceSearch" sheet ---
$wsAdvanceSearch = $wb.Worksheets.Item("AdvanceSearch")
# J2 was "ProductQA" -> now "Dell"
$wsAdvanceSearch.Range("J2").Value = "Dell"
# I4 was "QA" -> now "ang" (quote-prefixed text, like H4 above)
$wsAdvanceSearch.Range("I4").Value = "'ang"
# I6 was "Murali" -> now "Radhe"
$wsAdvanceSearch.Range("I6").Value = "Radhe"

# --- Update selections / active sheet to match the latest save state ---
$wsShowDateRange.Range("H5").Select()

$wsAdvanceSearch.Activate()
$wsAdvanceSearch.Range("I5").Select()
